$wb = $excel.ActiveWorkbook

# StatOutput sheet: the breed filter changed (Akita -> English Setter query),
# so the cached stat counts for number_of_files / number_of_sample are now 0.
$statOutput = $wb.Worksheets.Item("StatOutput")

# Setting .Value to a numeric-looking string normally makes Excel store it
# as a number; force text entry by temporarily marking the cell as Text
# format, then clear the format override afterwards so the cell keeps its
# original (default/General) style while the stored value stays textual.
$statOutput.Range("A2").NumberFormat = "@"
$statOutput.Range("A2").Value = "0"
$statOutput.Range("A2").ClearFormats()

$statOutput.Range("B2").NumberFormat = "@"
$statOutput.Range("B2").Value = "0"
$statOutput.Range("B2").ClearFormats()

# StatOutput_Message sheet: update the second Cypher query (row 18) so it
# filters on 'English Setter' (matching the first query) instead of 'Akita'.
$statMsg = $wb.Worksheets.Item("StatOutput_Message")
$newCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['English Setter']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$statMsg.Range("A18").Value = $newCypher
